$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to remain text so numeric-looking strings
# (e.g. "521.30", "1.00") are not coerced into numbers, losing formatting/precision.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '58.700.14'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '2.518.24'
$ws.Range('E3').Value = '  +2.59%  '
$ws.Range('D4').Value = '0.996'
$ws.Range('E4').Value = '  -0.39%  '
$ws.Range('D5').Value = '521.30'
$ws.Range('E5').Value = '  +0.93%  '
$ws.Range('D6').Value = '132.89'
$ws.Range('E6').Value = '  +0.55%  '
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('D8').Value = '0.557'
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('D9').Value = '2.517.47'
$ws.Range('E9').Value = '  +2.43%  '
$ws.Range('D10').Value = '0.0977'
$ws.Range('E10').Value = '  -0.42%  '
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '5.17'
$ws.Range('E12').Value = '  -1.65%  '
$ws.Range('D13').Value = '0.333'
$ws.Range('E13').Value = '  -1.91%  '
$ws.Range('D14').Value = '2.954.87'
$ws.Range('E14').Value = '  +2.16%  '
$ws.Range('D15').Value = '58.612.41'
$ws.Range('E15').Value = '  +1.32%  '
$ws.Range('D16').Value = '22.19'
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('E17').Value = '  +0.19%  '
$ws.Range('D18').Value = '2.503.09'
$ws.Range('E18').Value = '  +2.10%  '
$ws.Range('D19').Value = '10.72'
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('D20').Value = '321.62'
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('D21').Value = '4.16'
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('D22').Value = '6.03'
$ws.Range('E22').Value = '  +5.26%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').Value = '64.44'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('E25').Value = '  -1.71%  '
$ws.Range('E26').Value = '  +1.04%  '
$ws.Range('D27').Value = '0.990'
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('D28').Value = '7.39'
$ws.Range('E28').Value = '  +0.96%  '
$ws.Range('E29').Value = '  +2.38%  '
$ws.Range('E30').Value = '  +1.79%  '
$ws.Range('E31').Value = '  +3.48%  '
$ws.Range('D32').Value = '167.41'
$ws.Range('E32').Value = '  +0.53%  '
$ws.Range('D33').Value = '6.27'
$ws.Range('E33').Value = '  +1.18%  '
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Value = '18.10'
$ws.Range('E36').Value = '  +0.31%  '
$ws.Range('E37').Value = '  -2.06%  '
$ws.Range('D38').Value = '3.97'
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('D39').Value = '36.86'
$ws.Range('E39').Value = '  +1.85%  '
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').Value = '0.783'
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('D42').Value = '278.62'
$ws.Range('E42').Value = '  +2.97%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = '3.48'
$ws.Range('E43').Value = '  +1.57%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '5.10'
$ws.Range('E44').Value = '  +1.91%  '
$ws.Range('D45').Value = '0.600'
$ws.Range('E45').Value = '  +2.11%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').Value = '0.0918'
$ws.Range('E46').Value = '  +1.33%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '122.29'
$ws.Range('E47').Value = '  -1.99%  '
$ws.Range('D48').Value = '0.0502'
$ws.Range('E48').Value = '  +3.59%  '
$ws.Range('E49').Value = '  +1.08%  '
$ws.Range('D50').Value = '0.0214'
$ws.Range('E50').Value = '  +1.80%  '
$ws.Range('D51').Value = '16.94'
$ws.Range('E51').Value = '  +1.62%  '
